$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.214638590812683
$ws.Range("B1").Value = 4.64321756362915
$ws.Range("C1").Value = 4.410257339477539
$ws.Range("D1").Value = 1.985043406486511
$ws.Range("E1").Value = 1.287562489509583
